$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-02 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-03 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("734×7=5138", $true, $false, $false, $false, $false, $true, 1, $false, "789×4=3156", 2) | Out-Null
$d.Content.Find.Execute("743×6=4458", $true, $false, $false, $false, $false, $true, 1, $false, "362×3=1086", 2) | Out-Null
$d.Content.Find.Execute("343×5=1715", $true, $false, $false, $false, $false, $true, 1, $false, "791×4=3164", 2) | Out-Null
$d.Content.Find.Execute("623×5=3115", $true, $false, $false, $false, $false, $true, 1, $false, "201×9=1809", 2) | Out-Null
$d.Content.Find.Execute("558×7=3906", $true, $false, $false, $false, $false, $true, 1, $false, "329×8=2632", 2) | Out-Null
$d.Content.Find.Execute("401×4=1604", $true, $false, $false, $false, $false, $true, 1, $false, "262×7=1834", 2) | Out-Null
$d.Content.Find.Execute("528×7=3696", $true, $false, $false, $false, $false, $true, 1, $false, "225×7=1575", 2) | Out-Null
$d.Content.Find.Execute("457×8=3656", $true, $false, $false, $false, $false, $true, 1, $false, "443×3=1329", 2) | Out-Null
$d.Content.Find.Execute("329×7=2303", $true, $false, $false, $false, $false, $true, 1, $false, "867×6=5202", 2) | Out-Null
$d.Content.Find.Execute("234×6=1404", $true, $false, $false, $false, $false, $true, 1, $false, "553×3=1659", 2) | Out-Null
$d.Content.Find.Execute("301×5=1505", $true, $false, $false, $false, $false, $true, 1, $false, "281×3=843", 2) | Out-Null
$d.Content.Find.Execute("826×8=6608", $true, $false, $false, $false, $false, $true, 1, $false, "425×8=3400", 2) | Out-Null
$d.Content.Find.Execute("386×9=3474", $true, $false, $false, $false, $false, $true, 1, $false, "677×9=6093", 2) | Out-Null
$d.Content.Find.Execute("702×7=4914", $true, $false, $false, $false, $false, $true, 1, $false, "288×2=576", 2) | Out-Null
$d.Content.Find.Execute("448×4=1792", $true, $false, $false, $false, $false, $true, 1, $false, "413×3=1239", 2) | Out-Null
$d.Content.Find.Execute("713×6=4278", $true, $false, $false, $false, $false, $true, 1, $false, "847×7=5929", 2) | Out-Null
$d.Content.Find.Execute("868×6=5208", $true, $false, $false, $false, $false, $true, 1, $false, "166×6=996", 2) | Out-Null
$d.Content.Find.Execute("710×6=4260", $true, $false, $false, $false, $false, $true, 1, $false, "199×8=1592", 2) | Out-Null
$d.Content.Find.Execute("458×4=1832", $true, $false, $false, $false, $false, $true, 1, $false, "974×7=6818", 2) | Out-Null
$d.Content.Find.Execute("911×9=8199", $true, $false, $false, $false, $false, $true, 1, $false, "993×6=5958", 2) | Out-Null
$d.Content.Find.Execute("432×9=3888", $true, $false, $false, $false, $false, $true, 1, $false, "977×9=8793", 2) | Out-Null
$d.Content.Find.Execute("227×8=1816", $true, $false, $false, $false, $false, $true, 1, $false, "406×9=3654", 2) | Out-Null
$d.Content.Find.Execute("746×2=1492", $true, $false, $false, $false, $false, $true, 1, $false, "577×9=5193", 2) | Out-Null
$d.Content.Find.Execute("257×3=771", $true, $false, $false, $false, $false, $true, 1, $false, "239×4=956", 2) | Out-Null
$d.Content.Find.Execute("980×2=1960", $true, $false, $false, $false, $false, $true, 1, $false, "580×2=1160", 2) | Out-Null

Write-Host "Replacements complete"
